$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.567.40'
$ws.Range("E2").Value = '  +5.14%  '
$ws.Range("D3").Value = '3.638.19'
$ws.Range("E3").Value = '  +4.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '194.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.42%  '
$ws.Range("E7").Value = '  +2.21%  '
$ws.Range("D8").Value = '3.633.01'
$ws.Range("E8").Value = '  +4.88%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +4.91%  '
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000292'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.89%  '
$ws.Range("D15").Value = '4.209.10'
$ws.Range("E15").Value = '  +4.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.94%  '
$ws.Range("D17").Value = '3.629.58'
$ws.Range("E17").Value = '  +4.57%  '
$ws.Range("D18").Value = '70.474.70'
$ws.Range("E18").Value = '  +5.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.83%  '
$ws.Range("E20").Value = '  +2.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("E27").Value = '  +7.64%  '
$ws.Range("E28").Value = '  +5.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.05%  '
$ws.Range("E32").Value = '  +8.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '628.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '41.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.414'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.53%  '
$ws.Range("D38").Value = '0.0₃0828'
$ws.Range("E38").Value = '  +9.14%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("D42").Value = '3.306.56'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0455'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.76%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.139'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.25%  '
